$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blancos")

$ws.Cells.Item(2, 5).Value = 'FÍSICA I'
$ws.Cells.Item(2, 6).Value = 'Polanco Domínguez Rosa María'
$ws.Cells.Item(4, 5).Value = 'ECOLOGÍA'
$ws.Cells.Item(4, 6).Value = 'Camarillo Aburto Raymundo'
$ws.Cells.Item(5, 5).Value = 'INGLÉS IV'
$ws.Cells.Item(5, 6).Value = 'González Nuñez Veronica'
$ws.Cells.Item(6, 5).Value = 'CÁLCULO DIFERENCIAL'
$ws.Cells.Item(6, 6).Value = 'Ortega Valle Manuel'
$ws.Cells.Item(7, 5).Value = 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO'
$ws.Cells.Item(7, 6).Value = 'Zarate Amezcua Eladio Jorge'
$ws.Cells.Item(8, 5).Value = 'FÍSICA I'
$ws.Cells.Item(8, 6).Value = 'Polanco Domínguez Rosa María'
$ws.Cells.Item(9, 5).Value = 'ECOLOGÍA'
$ws.Cells.Item(9, 6).Value = 'Camarillo Aburto Raymundo'
$ws.Cells.Item(10, 5).Value = 'INGLÉS IV'
$ws.Cells.Item(10, 6).Value = 'González Nuñez Veronica'
$ws.Cells.Item(13, 5).Value = 'FÍSICA I'
$ws.Cells.Item(13, 6).Value = 'Polanco Domínguez Rosa María'
$ws.Cells.Item(15, 5).Value = 'ECOLOGÍA'
$ws.Cells.Item(15, 6).Value = 'Camarillo Aburto Raymundo'
$ws.Cells.Item(16, 5).Value = 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO'
$ws.Cells.Item(16, 6).Value = 'Zarate Amezcua Eladio Jorge'
$ws.Cells.Item(17, 5).Value = 'CÁLCULO DIFERENCIAL'
$ws.Cells.Item(17, 6).Value = 'Ortega Valle Manuel'
$ws.Cells.Item(18, 5).Value = 'ECOLOGÍA'
$ws.Cells.Item(18, 6).Value = 'Camarillo Aburto Raymundo'
$ws.Cells.Item(19, 5).Value = 'INGLÉS IV'
$ws.Cells.Item(19, 6).Value = 'González Nuñez Veronica'
$ws.Cells.Item(20, 5).Value = 'CÁLCULO DIFERENCIAL'
$ws.Cells.Item(20, 6).Value = 'Ortega Valle Manuel'
$ws.Cells.Item(21, 5).Value = 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO'
$ws.Cells.Item(21, 6).Value = 'Zarate Amezcua Eladio Jorge'
$ws.Cells.Item(22, 5).Value = 'FÍSICA I'
$ws.Cells.Item(22, 6).Value = 'Polanco Domínguez Rosa María'
$ws.Cells.Item(23, 5).Value = 'ECOLOGÍA'
$ws.Cells.Item(23, 6).Value = 'Camarillo Aburto Raymundo'
$ws.Cells.Item(24, 5).Value = 'CÁLCULO DIFERENCIAL'
$ws.Cells.Item(24, 6).Value = 'Ortega Valle Manuel'
$ws.Cells.Item(25, 5).Value = 'FÍSICA I'
$ws.Cells.Item(25, 6).Value = 'Polanco Domínguez Rosa María'
$ws.Cells.Item(29, 5).Value = 'FÍSICA I'
$ws.Cells.Item(29, 6).Value = 'Polanco Domínguez Rosa María'
$ws.Cells.Item(31, 5).Value = 'ECOLOGÍA'
$ws.Cells.Item(31, 6).Value = 'Camarillo Aburto Raymundo'
$ws.Cells.Item(39, 5).Value = 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO'
$ws.Cells.Item(39, 6).Value = 'Zarate Amezcua Eladio Jorge'
$ws.Cells.Item(40, 5).Value = 'CÁLCULO DIFERENCIAL'
$ws.Cells.Item(40, 6).Value = 'Ortega Valle Manuel'
$ws.Cells.Item(41, 5).Value = 'FÍSICA I'
$ws.Cells.Item(41, 6).Value = 'Polanco Domínguez Rosa María'
$ws.Cells.Item(42, 5).Value = 'INGLÉS IV'
$ws.Cells.Item(42, 6).Value = 'González Nuñez Veronica'
$ws.Cells.Item(43, 5).Value = 'ECOLOGÍA'
$ws.Cells.Item(43, 6).Value = 'Camarillo Aburto Raymundo'
$ws.Cells.Item(45, 5).Value = 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO'
$ws.Cells.Item(45, 6).Value = 'Zarate Amezcua Eladio Jorge'
$ws.Cells.Item(46, 5).Value = 'ECOLOGÍA'
$ws.Cells.Item(46, 6).Value = 'Camarillo Aburto Raymundo'
$ws.Cells.Item(47, 5).Value = 'FÍSICA I'
$ws.Cells.Item(47, 6).Value = 'Polanco Domínguez Rosa María'
$ws.Cells.Item(48, 5).Value = 'INGLÉS IV'
$ws.Cells.Item(48, 6).Value = 'González Nuñez Veronica'
$ws.Cells.Item(50, 5).Value = 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO'
$ws.Cells.Item(50, 6).Value = 'Zarate Amezcua Eladio Jorge'
$ws.Cells.Item(51, 5).Value = 'ECOLOGÍA'
$ws.Cells.Item(51, 6).Value = 'Camarillo Aburto Raymundo'
$ws.Cells.Item(52, 5).Value = 'INGLÉS IV'
$ws.Cells.Item(52, 6).Value = 'González Nuñez Veronica'
$ws.Cells.Item(53, 5).Value = 'FÍSICA I'
$ws.Cells.Item(53, 6).Value = 'Polanco Domínguez Rosa María'
$ws.Cells.Item(55, 5).Value = 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO'
$ws.Cells.Item(55, 6).Value = 'Zarate Amezcua Eladio Jorge'
$ws.Cells.Item(56, 5).Value = 'INGLÉS IV'
$ws.Cells.Item(56, 6).Value = 'González Nuñez Veronica'
$ws.Cells.Item(57, 5).Value = 'ECOLOGÍA'
$ws.Cells.Item(57, 6).Value = 'Camarillo Aburto Raymundo'
$ws.Cells.Item(58, 5).Value = 'CÁLCULO DIFERENCIAL'
$ws.Cells.Item(58, 6).Value = 'Ortega Valle Manuel'
$ws.Cells.Item(59, 5).Value = 'INGLÉS IV'
$ws.Cells.Item(59, 6).Value = 'González Nuñez Veronica'
$ws.Cells.Item(60, 5).Value = 'ECOLOGÍA'
$ws.Cells.Item(60, 6).Value = 'Camarillo Aburto Raymundo'
$ws.Cells.Item(61, 5).Value = 'FÍSICA I'
$ws.Cells.Item(61, 6).Value = 'Polanco Domínguez Rosa María'
$ws.Cells.Item(62, 5).Value = 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO'
$ws.Cells.Item(62, 6).Value = 'Zarate Amezcua Eladio Jorge'
$ws.Cells.Item(64, 5).Value = 'CÁLCULO DIFERENCIAL'
$ws.Cells.Item(64, 6).Value = 'Ortega Valle Manuel'
$ws.Cells.Item(65, 5).Value = 'ECOLOGÍA'
$ws.Cells.Item(65, 6).Value = 'Camarillo Aburto Raymundo'
$ws.Cells.Item(66, 5).Value = 'FÍSICA I'
$ws.Cells.Item(66, 6).Value = 'Polanco Domínguez Rosa María'
$ws.Cells.Item(67, 5).Value = 'FÍSICA I'
$ws.Cells.Item(67, 6).Value = 'Polanco Domínguez Rosa María'
$ws.Cells.Item(69, 5).Value = 'ECOLOGÍA'
$ws.Cells.Item(69, 6).Value = 'Camarillo Aburto Raymundo'
$ws.Cells.Item(70, 5).Value = 'CÁLCULO DIFERENCIAL'
$ws.Cells.Item(70, 6).Value = 'Ortega Valle Manuel'
$ws.Cells.Item(71, 5).Value = 'FÍSICA I'
$ws.Cells.Item(71, 6).Value = 'Polanco Domínguez Rosa María'
$ws.Cells.Item(72, 5).Value = 'INGLÉS IV'
$ws.Cells.Item(72, 6).Value = 'González Nuñez Veronica'
$ws.Cells.Item(73, 5).Value = 'ECOLOGÍA'
$ws.Cells.Item(73, 6).Value = 'Camarillo Aburto Raymundo'
$ws.Cells.Item(74, 5).Value = 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO'
$ws.Cells.Item(74, 6).Value = 'Zarate Amezcua Eladio Jorge'
$ws.Cells.Item(76, 5).Value = 'FÍSICA I'
$ws.Cells.Item(76, 6).Value = 'Polanco Domínguez Rosa María'
$ws.Cells.Item(77, 5).Value = 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO'
$ws.Cells.Item(77, 6).Value = 'Zarate Amezcua Eladio Jorge'
$ws.Cells.Item(78, 5).Value = 'INGLÉS IV'
$ws.Cells.Item(78, 6).Value = 'González Nuñez Veronica'
$ws.Cells.Item(80, 5).Value = 'ECOLOGÍA'
$ws.Cells.Item(80, 6).Value = 'Camarillo Aburto Raymundo'
$ws.Cells.Item(81, 5).Value = 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO'
$ws.Cells.Item(81, 6).Value = 'Zarate Amezcua Eladio Jorge'
$ws.Cells.Item(82, 5).Value = 'FÍSICA I'
$ws.Cells.Item(82, 6).Value = 'Polanco Domínguez Rosa María'
$ws.Cells.Item(83, 5).Value = 'CÁLCULO DIFERENCIAL'
$ws.Cells.Item(83, 6).Value = 'Ortega Valle Manuel'
$ws.Cells.Item(84, 5).Value = 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO'
$ws.Cells.Item(84, 6).Value = 'Zarate Amezcua Eladio Jorge'
$ws.Cells.Item(86, 5).Value = 'ECOLOGÍA'
$ws.Cells.Item(86, 6).Value = 'Camarillo Aburto Raymundo'
$ws.Cells.Item(87, 5).Value = 'INGLÉS IV'
$ws.Cells.Item(87, 6).Value = 'González Nuñez Veronica'
$ws.Cells.Item(90, 5).Value = 'FÍSICA I'
$ws.Cells.Item(90, 6).Value = 'Polanco Domínguez Rosa María'
$ws.Cells.Item(91, 5).Value = 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO'
$ws.Cells.Item(91, 6).Value = 'Zarate Amezcua Eladio Jorge'
$ws.Cells.Item(92, 5).Value = 'INGLÉS IV'
$ws.Cells.Item(92, 6).Value = 'González Nuñez Veronica'
$ws.Cells.Item(93, 5).Value = 'ECOLOGÍA'
$ws.Cells.Item(93, 6).Value = 'Camarillo Aburto Raymundo'
$ws.Cells.Item(94, 5).Value = 'CÁLCULO DIFERENCIAL'
$ws.Cells.Item(94, 6).Value = 'Ortega Valle Manuel'
$ws.Cells.Item(95, 5).Value = 'CÁLCULO DIFERENCIAL'
$ws.Cells.Item(95, 6).Value = 'Ortega Valle Manuel'
$ws.Cells.Item(96, 5).Value = 'ECOLOGÍA'
$ws.Cells.Item(96, 6).Value = 'Camarillo Aburto Raymundo'
$ws.Cells.Item(97, 5).Value = 'FÍSICA I'
$ws.Cells.Item(97, 6).Value = 'Polanco Domínguez Rosa María'
$ws.Cells.Item(98, 5).Value = 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO'
$ws.Cells.Item(98, 6).Value = 'Zarate Amezcua Eladio Jorge'
$ws.Cells.Item(99, 5).Value = 'INGLÉS IV'
$ws.Cells.Item(99, 6).Value = 'González Nuñez Veronica'
$ws.Cells.Item(100, 5).Value = 'ECOLOGÍA'
$ws.Cells.Item(100, 6).Value = 'Camarillo Aburto Raymundo'
$ws.Cells.Item(102, 5).Value = 'CÁLCULO DIFERENCIAL'
$ws.Cells.Item(102, 6).Value = 'Ortega Valle Manuel'
$ws.Cells.Item(103, 5).Value = 'MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO'
$ws.Cells.Item(103, 6).Value = 'Zarate Amezcua Eladio Jorge'
